$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = 42
$ws.Range("E6").Value = 46

$ws.Range("M21").Select()
